$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 318
$ws1.Range("F8").Value = 2083
$ws1.Range("F10").Value = 46
$ws1.Range("F11").Value = 44
$ws1.Range("F14").Value = 1345
$ws1.Range("F15").Value = 60
$ws1.Range("F19").Value = 531
$ws1.Range("F21").Value = 151
$ws1.Range("F22").Value = 7147
$ws1.Range("F23").Value = 7786
$ws1.Range("F32").Value = 14
$ws1.Range("F36").Value = 1410
$ws1.Range("F40").Value = 286
$ws1.Range("F41").Value = 708
$ws1.Range("F43").Value = 1358
$ws1.Range("F45").Value = 238
$ws1.Range("F48").Value = 154
$ws1.Range("F49").Value = 154

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 31

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 181
$ws3.Range("F3").Value = 2604
$ws3.Range("F4").Value = 277
$ws3.Range("F5").Value = 130

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 277
$ws4.Range("F5").Value = 130
$ws4.Range("F7").Value = 318
$ws4.Range("F8").Value = 31
$ws4.Range("F10").Value = 2083
$ws4.Range("F11").Value = 46
$ws4.Range("F12").Value = 44
$ws4.Range("F16").Value = 1345
$ws4.Range("F17").Value = 60
$ws4.Range("F20").Value = 531
$ws4.Range("F22").Value = 7147
$ws4.Range("F23").Value = 7786
$ws4.Range("F29").Value = 14
$ws4.Range("F31").Value = 1410
$ws4.Range("F35").Value = 286
$ws4.Range("F38").Value = 708
$ws4.Range("F43").Value = 238
$ws4.Range("F46").Value = 154
$ws4.Range("F47").Value = 154
